$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update rows 2-31, columns A-F, per the new stock-symbol lists
$ws.Cells.Item(2, 1).Value = 0
$ws.Cells.Item(2, 2).Value = "NSE:ALICON"
$ws.Cells.Item(2, 3).Value = "NSE:APARINDS"
$ws.Cells.Item(2, 4).Value = ""
$ws.Cells.Item(2, 5).Value = ""
$ws.Cells.Item(2, 6).Value = "NSE:BEL"

$ws.Cells.Item(3, 1).Value = 1
$ws.Cells.Item(3, 2).Value = "NSE:ALLCARGO"
$ws.Cells.Item(3, 3).Value = "NSE:GSS"
$ws.Cells.Item(3, 4).Value = ""
$ws.Cells.Item(3, 5).Value = ""
$ws.Cells.Item(3, 6).Value = "NSE:IRCTC"

$ws.Cells.Item(4, 1).Value = 2
$ws.Cells.Item(4, 2).Value = "NSE:AMBER"
$ws.Cells.Item(4, 3).Value = "NSE:KAPSTON"
$ws.Cells.Item(4, 4).Value = ""
$ws.Cells.Item(4, 5).Value = ""
$ws.Cells.Item(4, 6).Value = ""

$ws.Cells.Item(5, 1).Value = 3
$ws.Cells.Item(5, 2).Value = "NSE:BDL"
$ws.Cells.Item(5, 3).Value = ""
$ws.Cells.Item(5, 4).Value = ""
$ws.Cells.Item(5, 5).Value = ""
$ws.Cells.Item(5, 6).Value = ""

$ws.Cells.Item(6, 1).Value = 4
$ws.Cells.Item(6, 2).Value = "NSE:BEDMUTHA"
$ws.Cells.Item(6, 3).Value = ""
$ws.Cells.Item(6, 4).Value = ""
$ws.Cells.Item(6, 5).Value = ""
$ws.Cells.Item(6, 6).Value = ""

$ws.Cells.Item(7, 1).Value = 5
$ws.Cells.Item(7, 2).Value = "NSE:BEL"
$ws.Cells.Item(7, 3).Value = ""
$ws.Cells.Item(7, 4).Value = ""
$ws.Cells.Item(7, 5).Value = ""
$ws.Cells.Item(7, 6).Value = ""

$ws.Cells.Item(8, 1).Value = 6
$ws.Cells.Item(8, 2).Value = "NSE:BHEL"
$ws.Cells.Item(8, 3).Value = ""
$ws.Cells.Item(8, 4).Value = ""
$ws.Cells.Item(8, 5).Value = ""
$ws.Cells.Item(8, 6).Value = ""

$ws.Cells.Item(9, 1).Value = 7
$ws.Cells.Item(9, 2).Value = "NSE:BIOFILCHEM"
$ws.Cells.Item(9, 3).Value = ""
$ws.Cells.Item(9, 4).Value = ""
$ws.Cells.Item(9, 5).Value = ""
$ws.Cells.Item(9, 6).Value = ""

$ws.Cells.Item(10, 1).Value = 8
$ws.Cells.Item(10, 2).Value = "NSE:CANTABIL"
$ws.Cells.Item(10, 3).Value = ""
$ws.Cells.Item(10, 4).Value = ""
$ws.Cells.Item(10, 5).Value = ""
$ws.Cells.Item(10, 6).Value = ""

$ws.Cells.Item(11, 1).Value = 9
$ws.Cells.Item(11, 2).Value = "NSE:DCMSRIND"
$ws.Cells.Item(11, 3).Value = ""
$ws.Cells.Item(11, 4).Value = ""
$ws.Cells.Item(11, 5).Value = ""
$ws.Cells.Item(11, 6).Value = ""

$ws.Cells.Item(12, 1).Value = 10
$ws.Cells.Item(12, 2).Value = "NSE:GEPIL"
$ws.Cells.Item(12, 3).Value = ""
$ws.Cells.Item(12, 4).Value = ""
$ws.Cells.Item(12, 5).Value = ""
$ws.Cells.Item(12, 6).Value = ""

$ws.Cells.Item(13, 1).Value = 11
$ws.Cells.Item(13, 2).Value = "NSE:HCC"
$ws.Cells.Item(13, 3).Value = ""
$ws.Cells.Item(13, 4).Value = ""
$ws.Cells.Item(13, 5).Value = ""
$ws.Cells.Item(13, 6).Value = ""

$ws.Cells.Item(14, 1).Value = 12
$ws.Cells.Item(14, 2).Value = "NSE:HEALTHY"
$ws.Cells.Item(14, 3).Value = ""
$ws.Cells.Item(14, 4).Value = ""
$ws.Cells.Item(14, 5).Value = ""
$ws.Cells.Item(14, 6).Value = ""

$ws.Cells.Item(15, 1).Value = 13
$ws.Cells.Item(15, 2).Value = "NSE:IFCI"
$ws.Cells.Item(15, 3).Value = ""
$ws.Cells.Item(15, 4).Value = ""
$ws.Cells.Item(15, 5).Value = ""
$ws.Cells.Item(15, 6).Value = ""

$ws.Cells.Item(16, 1).Value = 14
$ws.Cells.Item(16, 2).Value = "NSE:IFGLEXPOR"
$ws.Cells.Item(16, 3).Value = ""
$ws.Cells.Item(16, 4).Value = ""
$ws.Cells.Item(16, 5).Value = ""
$ws.Cells.Item(16, 6).Value = ""

$ws.Cells.Item(17, 1).Value = 15
$ws.Cells.Item(17, 2).Value = "NSE:INDIANHUME"
$ws.Cells.Item(17, 3).Value = ""
$ws.Cells.Item(17, 4).Value = ""
$ws.Cells.Item(17, 5).Value = ""
$ws.Cells.Item(17, 6).Value = ""

$ws.Cells.Item(18, 1).Value = 16
$ws.Cells.Item(18, 2).Value = "NSE:IRCTC"
$ws.Cells.Item(18, 3).Value = ""
$ws.Cells.Item(18, 4).Value = ""
$ws.Cells.Item(18, 5).Value = ""
$ws.Cells.Item(18, 6).Value = ""

$ws.Cells.Item(19, 1).Value = 17
$ws.Cells.Item(19, 2).Value = "NSE:KOTAKBANK"
$ws.Cells.Item(19, 3).Value = ""
$ws.Cells.Item(19, 4).Value = ""
$ws.Cells.Item(19, 5).Value = ""
$ws.Cells.Item(19, 6).Value = ""

$ws.Cells.Item(20, 1).Value = 18
$ws.Cells.Item(20, 2).Value = "NSE:M&M"
$ws.Cells.Item(20, 3).Value = ""
$ws.Cells.Item(20, 4).Value = ""
$ws.Cells.Item(20, 5).Value = ""
$ws.Cells.Item(20, 6).Value = ""

$ws.Cells.Item(21, 1).Value = 19
$ws.Cells.Item(21, 2).Value = "NSE:MAITHANALL"
$ws.Cells.Item(21, 3).Value = ""
$ws.Cells.Item(21, 4).Value = ""
$ws.Cells.Item(21, 5).Value = ""
$ws.Cells.Item(21, 6).Value = ""

$ws.Cells.Item(22, 1).Value = 20
$ws.Cells.Item(22, 2).Value = "NSE:MSTCLTD"
$ws.Cells.Item(22, 3).Value = ""
$ws.Cells.Item(22, 4).Value = ""
$ws.Cells.Item(22, 5).Value = ""
$ws.Cells.Item(22, 6).Value = ""

$ws.Cells.Item(23, 1).Value = 21
$ws.Cells.Item(23, 2).Value = "NSE:NBIFIN"
$ws.Cells.Item(23, 3).Value = ""
$ws.Cells.Item(23, 4).Value = ""
$ws.Cells.Item(23, 5).Value = ""
$ws.Cells.Item(23, 6).Value = ""

$ws.Cells.Item(24, 1).Value = 22
$ws.Cells.Item(24, 2).Value = "NSE:NELCO"
$ws.Cells.Item(24, 3).Value = ""
$ws.Cells.Item(24, 4).Value = ""
$ws.Cells.Item(24, 5).Value = ""
$ws.Cells.Item(24, 6).Value = ""

$ws.Cells.Item(25, 1).Value = 23
$ws.Cells.Item(25, 2).Value = "NSE:NH"
$ws.Cells.Item(25, 3).Value = ""
$ws.Cells.Item(25, 4).Value = ""
$ws.Cells.Item(25, 5).Value = ""
$ws.Cells.Item(25, 6).Value = ""

$ws.Cells.Item(26, 1).Value = 24
$ws.Cells.Item(26, 2).Value = "NSE:NYKAA"
$ws.Cells.Item(26, 3).Value = ""
$ws.Cells.Item(26, 4).Value = ""
$ws.Cells.Item(26, 5).Value = ""
$ws.Cells.Item(26, 6).Value = ""

$ws.Cells.Item(27, 1).Value = 25
$ws.Cells.Item(27, 2).Value = "NSE:PRINCEPIPE"
$ws.Cells.Item(27, 3).Value = ""
$ws.Cells.Item(27, 4).Value = ""
$ws.Cells.Item(27, 5).Value = ""
$ws.Cells.Item(27, 6).Value = ""

$ws.Cells.Item(28, 1).Value = 26
$ws.Cells.Item(28, 2).Value = "NSE:RRKABEL"
$ws.Cells.Item(28, 3).Value = ""
$ws.Cells.Item(28, 4).Value = ""
$ws.Cells.Item(28, 5).Value = ""
$ws.Cells.Item(28, 6).Value = ""

$ws.Cells.Item(29, 1).Value = 27
$ws.Cells.Item(29, 2).Value = "NSE:RSYSTEMS"
$ws.Cells.Item(29, 3).Value = ""
$ws.Cells.Item(29, 4).Value = ""
$ws.Cells.Item(29, 5).Value = ""
$ws.Cells.Item(29, 6).Value = ""

$ws.Cells.Item(30, 1).Value = 28
$ws.Cells.Item(30, 2).Value = "NSE:RVNL"
$ws.Cells.Item(30, 3).Value = ""
$ws.Cells.Item(30, 4).Value = ""
$ws.Cells.Item(30, 5).Value = ""
$ws.Cells.Item(30, 6).Value = ""

$ws.Cells.Item(31, 1).Value = 29
$ws.Cells.Item(31, 2).Value = "NSE:SAKSOFT"
$ws.Cells.Item(31, 3).Value = ""
$ws.Cells.Item(31, 4).Value = ""
$ws.Cells.Item(31, 5).Value = ""
$ws.Cells.Item(31, 6).Value = ""

